$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) DataSet sheet: add new cell X7 with URL text + hyperlink,
#    using the existing Hyperlink-ish style already used by G2/G7.
#    (This also creates shared string #179, the URL text, FIRST so
#    the shared-string table ordering lines up with the target file.)
# ------------------------------------------------------------------
$wsDataSet = $wb.Worksheets.Item("DataSet")

$hydroUrl = "https://mcloud-na-preprod.hydroflask.com/hydroqalotuswvae"
$wsDataSet.Range("X7").Value = $hydroUrl
$wsDataSet.Hyperlinks.Add($wsDataSet.Range("X7"), $hydroUrl) | Out-Null
$wsDataSet.Range("X7").Style = "Hyperlink"

# ------------------------------------------------------------------
# 2) Forms sheet: remove the Products / Quantity / DOB columns
#    (N:P), repurpose the old "Confirm Password" header cell as
#    "HydroAnswers", and append a new data row describing the chat
#    options / categories.
# ------------------------------------------------------------------
$wsForms = $wb.Worksheets.Item("Forms")

# Record the hyperlinks that exist before the column shuffle so we
# can recreate them at their correct post-delete addresses (the
# column delete does not itself repoint hyperlink ranges).
$formsLinks = @()
foreach ($hl in $wsForms.Hyperlinks) {
    $formsLinks += , @($hl.Range.Address(), $hl.Address)
}

$wsForms.Columns("N:P").Delete()

$wsForms.Cells.Hyperlinks.Delete()
$wsForms.Hyperlinks.Add($wsForms.Range("G2"), "mailto:qatesting.lotuswave@gmail.com") | Out-Null
$wsForms.Range("G2").Style = "Hyperlink"
$wsForms.Hyperlinks.Add($wsForms.Range("Q2"), "http://www.lotuswavess.com/") | Out-Null
$wsForms.Range("Q2").Style = "Hyperlink"
$wsForms.Hyperlinks.Add($wsForms.Range("G3"), "mailto:qatesting.lotuswave@gmail.com") | Out-Null
$wsForms.Range("G3").Style = "Hyperlink"

$wsForms.Range("A4").Value = "Chatoptions"
$wsForms.Range("D1").Value = "HydroAnswers"
$wsForms.Range("D4").Value = "General,Order Inquiries,Product Info,Warranty,Retail Information"

# ------------------------------------------------------------------
# 3) Walk the sheets the way the author apparently did: click around
#    on each tab (updating the stored selection on every sheet) and
#    finish with "Forms" as the active tab.
# ------------------------------------------------------------------
$wsDataSet.Activate()
$wsDataSet.Range("A17").Select()

$wsSheet1 = $wb.Worksheets.Item("Sheet1")
$wsSheet1.Activate()

$wsPDP = $wb.Worksheets.Item("PDP")
$wsPDP.Activate()
$wsPDP.Range("J11").Select()

$wsBundle = $wb.Worksheets.Item("Bundle")
$wsBundle.Activate()
$wsBundle.Range("A6").Select()

$wsForms.Activate()
$wsForms.Range("F14").Select()
